# Refresh cryptos list: updated prices and 1h volume-change percentages
# (and corrected the Hedera/VeChain row ordering) from the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (e.g. "27.733.23"); force the
# whole column to Text format first so Excel does not reinterpret the
# plain-decimal ones (e.g. "327.29") as numbers, then drop the explicit
# format again so the cells end up with the same (default) style as before.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.733.23"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.774.54"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "327.29"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.4582"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").Value = "0.3583"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").Value = "0.07493"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "1.104"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "6.044"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "7.226"
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "1.778.15"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "93.77"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "0.06435"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "17.09"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").Value = "5.810"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("D23").Value = "27.798.10"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "164.57"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").Value = "20.27"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "1.978.93"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D29").Value = "2.167"
$ws.Range("E29").Value = "  +4.00%  "
$ws.Range("D30").Value = "125.79"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").Value = "0.09227"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "3.670"
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("D34").Value = "5.535"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "11.88"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "0.02298"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06197"
$ws.Range("E37").Value = "  +3.04%  "
$ws.Range("D38").Value = "0.2091"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "4.961"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "1.190"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").Value = "1.390"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "7.801"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "13.24"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "3.745"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "0.5904"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "122.61"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "1.952"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "0.06925"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").Value = "1.140"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "72.39"
$ws.Range("E51").Value = "  +0.61%  "

$priceRange.ClearFormats()
